# "A bit more stripping on uneeded styles from menu bar"
#
# Turn off the "Quick Style" (Styles Gallery) flag -- i.e. remove the
# w:qFormat marker -- for a further batch of built-in styles. This is
# the same kind of cleanup as an earlier pass, just covering more
# styles: the heading 1-3 styles, Title, Subtitle, Quote, Intense
# Emphasis, Intense Quote and Intense Reference should no longer be
# promoted as "quick styles" in the Word Styles gallery / menu bar.

$d = $word.ActiveDocument

$styleNames = @(
    "Heading 1",
    "Heading 2",
    "Heading 3",
    "Title",
    "Subtitle",
    "Quote",
    "Intense Emphasis",
    "Intense Quote",
    "Intense Reference"
)

foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    if ($style -ne $null) {
        $style.QuickStyle = $false
    }
}
